$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1228.3529
$ws.Range("I9").Value = 1353.6666
$ws.Range("J9").Value = 1087.375
$ws.Range("K9").Value = 1353.6666
$ws.Range("L9").Value = 1087.375
$ws.Range("M9").Value = -1184.6666
$ws.Range("N9").Value = -1425.375
$ws.Range("H53").Value = 6687.5264
$ws.Range("J53").Value = 14523.375
$ws.Range("L53").Value = 14523.375
$ws.Range("N53").Value = -15797.375
$ws.Range("H64").Value = 4600
$ws.Range("I64").Value = 4000
$ws.Range("J64").Value = 4800
$ws.Range("K64").Value = 4000
$ws.Range("L64").Value = 4800
$ws.Range("M64").Value = -3752
$ws.Range("N64").Value = -5296
$ws.Range("H67").Value = 4600
$ws.Range("I67").Value = 4000
$ws.Range("J67").Value = 4800
$ws.Range("K67").Value = 4000
$ws.Range("L67").Value = 4800
$ws.Range("M67").Value = -3142
$ws.Range("N67").Value = -6516
$ws.Range("H100").Value = 3223.7896
$ws.Range("I100").Value = 1796.6364
$ws.Range("J100").Value = 5186.125
$ws.Range("K100").Value = 1796.6364
$ws.Range("L100").Value = 5186.125
$ws.Range("M100").Value = -1255.6364
$ws.Range("N100").Value = -6268.125
$ws.Range("H125").Value = 4135945.8
$ws.Range("I125").Value = 5054856
$ws.Range("J125").Value = 849.5
$ws.Range("K125").Value = 45493704
$ws.Range("L125").Value = 7645.5
$ws.Range("M125").Value = -45491244
$ws.Range("N125").Value = -12565.5
$ws.Range("H132").Value = 17684.79
$ws.Range("I132").Value = 13588.529
$ws.Range("J132").Value = 52503
$ws.Range("K132").Value = 40765.587
$ws.Range("L132").Value = 157509
$ws.Range("M132").Value = -38235.587
$ws.Range("N132").Value = -162569
$ws.Range("H137").Value = 9837.325000000001
$ws.Range("I137").Value = 2742.261
$ws.Range("J137").Value = 17996.65
$ws.Range("K137").Value = 8226.782999999999
$ws.Range("L137").Value = 53989.95
$ws.Range("M137").Value = -5676.782999999999
$ws.Range("N137").Value = -59089.95
$ws.Range("H138").Value = 3668.5
$ws.Range("I138").Value = 3265.7273
$ws.Range("J138").Value = 4554.6
$ws.Range("K138").Value = 9797.1819
$ws.Range("L138").Value = 13663.8
$ws.Range("M138").Value = -4657.1819
$ws.Range("N138").Value = -23943.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 6666
$ws.Range("J34").Value = 6666
$ws.Range("L34").Value = 6666
$ws.Range("M34").Value = -7208
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H61").Value = 23578.766
$ws.Range("I61").Value = 3359.4
$ws.Range("K61").Value = 3359.4
$ws.Range("M61").Value = -3147.4
$ws.Range("H74").Value = 23214
$ws.Range("I74").Value = 3802.4
$ws.Range("J74").Value = 42625.6
$ws.Range("K74").Value = 3802.4
$ws.Range("L74").Value = 42625.6
$ws.Range("M74").Value = -2928.4
$ws.Range("N74").Value = -44373.6
$ws.Range("H77").Value = 23214
$ws.Range("I77").Value = 3802.4
$ws.Range("J77").Value = 42625.6
$ws.Range("K77").Value = 19012
$ws.Range("L77").Value = 213128
$ws.Range("M77").Value = -14644
$ws.Range("N77").Value = -221864
$ws.Range("H102").Value = 20423.77
$ws.Range("I102").Value = 2215.7144
$ws.Range("K102").Value = 2215.7144
$ws.Range("M102").Value = -593.7143999999998
$ws.Range("H114").Value = 100398
$ws.Range("J114").Value = 100398
$ws.Range("L114").Value = 100398
$ws.Range("N114").Value = -109076
$ws.Range("H132").Value = 2870192.5
$ws.Range("I132").Value = 4585.1665
$ws.Range("J132").Value = 9122427
$ws.Range("K132").Value = 13755.4995
$ws.Range("L132").Value = 27367281
$ws.Range("M132").Value = -11225.4995
$ws.Range("N132").Value = -27372341
$ws.Range("H136").Value = 23578.766
$ws.Range("I136").Value = 3359.4
$ws.Range("K136").Value = 10078.2
$ws.Range("M136").Value = -7528.200000000001
$ws.Range("H141").Value = 35388
$ws.Range("I141").Value = 35388
$ws.Range("K141").Value = 35388
$ws.Range("M141").Value = -30208

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H44").Value = 40399.5
$ws.Range("J44").Value = 40399.5
$ws.Range("L44").Value = 40399.5
$ws.Range("N44").Value = -41393.5
$ws.Range("H99").Value = 7134.579
$ws.Range("I99").Value = 1334
$ws.Range("K99").Value = 1334
$ws.Range("M99").Value = 164
$ws.Range("H134").Value = 43165.6
$ws.Range("I134").Value = 5266.6665
$ws.Range("K134").Value = 15799.9995
$ws.Range("M134").Value = -13264.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 6500
$ws.Range("I62").Value = 10000
$ws.Range("K62").Value = 10000
$ws.Range("M62").Value = -9376
$ws.Range("H65").Value = 6500
$ws.Range("I65").Value = 10000
$ws.Range("K65").Value = 50000
$ws.Range("M65").Value = -46880

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 852.34784
$ws.Range("J23").Value = 1088.8334
$ws.Range("L23").Value = 3266.5002
$ws.Range("N23").Value = -3736.5002
$ws.Range("H74").Value = 16777.777
$ws.Range("J74").Value = 20714.285
$ws.Range("L74").Value = 62142.855
$ws.Range("N74").Value = -64264.855
$ws.Range("H77").Value = 16777.777
$ws.Range("J77").Value = 20714.285
$ws.Range("L77").Value = 186428.565
$ws.Range("N77").Value = -197036.565
$ws.Range("H106").Value = 4000
$ws.Range("J106").Value = 5666.6665
$ws.Range("L106").Value = 16999.9995
$ws.Range("N106").Value = -18891.9995
$ws.Range("H122").Value = 9359638
$ws.Range("J122").Value = 2580897.5
$ws.Range("L122").Value = 23228077.5
$ws.Range("N122").Value = -23232977.5
$ws.Range("H131").Value = 1456.33
$ws.Range("J131").Value = 1467.3334
$ws.Range("L131").Value = 4402.0002
$ws.Range("N131").Value = -14482.0002
$ws.Range("H140").Value = 1701.5
$ws.Range("I140").Value = 1701.5
$ws.Range("K140").Value = 5104.5
$ws.Range("M140").Value = 75.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 125093.94
$ws.Range("I12").Value = 125093.94
$ws.Range("K12").Value = 125093.94
$ws.Range("M12").Value = -124953.94
$ws.Range("H80").Value = 24008.389
$ws.Range("I80").Value = 20304.818
$ws.Range("J80").Value = 29828.285
$ws.Range("K80").Value = 20304.818
$ws.Range("L80").Value = 29828.285
$ws.Range("M80").Value = -19306.818
$ws.Range("N80").Value = -31824.285
$ws.Range("H83").Value = 24008.389
$ws.Range("I83").Value = 20304.818
$ws.Range("J83").Value = 29828.285
$ws.Range("K83").Value = 101524.09
$ws.Range("L83").Value = 149141.425
$ws.Range("M83").Value = -96532.09
$ws.Range("N83").Value = -159125.425
$ws.Range("H97").Value = 5857.8096
$ws.Range("I97").Value = 1554.3334
$ws.Range("K97").Value = 1554.3334
$ws.Range("M97").Value = -1058.3334
$ws.Range("H128").Value = 20000
$ws.Range("J128").Value = 20000
$ws.Range("L128").Value = 20000
$ws.Range("N128").Value = -29960
$ws.Range("H132").Value = 14396
$ws.Range("I132").Value = 8770
$ws.Range("J132").Value = 26585.666
$ws.Range("K132").Value = 26310
$ws.Range("L132").Value = 79756.99800000001
$ws.Range("M132").Value = -23780
$ws.Range("N132").Value = -84816.99800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 43483150
$ws.Range("I22").Value = 2736.75
$ws.Range("J22").Value = 142866960
$ws.Range("K22").Value = 2736.75
$ws.Range("L22").Value = 142866960
$ws.Range("M22").Value = -2441.75
$ws.Range("N22").Value = -142867550
$ws.Range("H25").Value = 245499.69
$ws.Range("H27").Value = 43483150
$ws.Range("I27").Value = 2736.75
$ws.Range("J27").Value = 142866960
$ws.Range("K27").Value = 2736.75
$ws.Range("L27").Value = 142866960
$ws.Range("M27").Value = -2629.75
$ws.Range("N27").Value = -142867174
$ws.Range("H38").Value = 85333.336
$ws.Range("I38").Value = 28000
$ws.Range("K38").Value = 28000
$ws.Range("M38").Value = -27590
$ws.Range("H46").Value = 3177.3125
$ws.Range("I46").Value = 1807.8334
$ws.Range("K46").Value = 1807.8334
$ws.Range("M46").Value = -1619.8334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 70007
$ws.Range("J45").Value = 70007
$ws.Range("L45").Value = 70007
$ws.Range("N45").Value = -70989
$ws.Range("H132").Value = 11978
$ws.Range("I132").Value = 4541
$ws.Range("K132").Value = 13623
$ws.Range("M132").Value = -11093
